$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("translation")

$ws.Range("A13").Value = "My Network"
$ws.Range("B13").Value = "My Network"
$ws.Range("C13").Value = "My Network"
$ws.Range("D13").Value = "My Network"
